# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# Values that are plain single-dot decimals (e.g. "241.40") are set through a
# temporary "@" (text) number format so Excel keeps them as literal text instead
# of silently parsing them into floating point numbers; the original cell style
# is restored immediately afterwards so no visible formatting changes.
$updates = @(
    @{ Cell = "D2"; Value = "29.165.06" }
    @{ Cell = "E2"; Value = "  +0.05%  " }
    @{ Cell = "D3"; Value = "1.833.99" }
    @{ Cell = "E3"; Value = "  -0.24%  " }
    @{ Cell = "D4"; Value = "0.9993" }
    @{ Cell = "E4"; Value = "  -0.03%  " }
    @{ Cell = "D5"; Value = "241.40" }
    @{ Cell = "E5"; Value = "  +0.60%  " }
    @{ Cell = "D6"; Value = "0.6647" }
    @{ Cell = "E6"; Value = "  -2.50%  " }
    @{ Cell = "D7"; Value = "1.0000" }
    @{ Cell = "E7"; Value = "  -0.01%  " }
    @{ Cell = "D8"; Value = "0.07421" }
    @{ Cell = "E8"; Value = "  -0.53%  " }
    @{ Cell = "D9"; Value = "0.2937" }
    @{ Cell = "E9"; Value = "  -2.00%  " }
    @{ Cell = "D10"; Value = "22.65" }
    @{ Cell = "E10"; Value = "  -2.41%  " }
    @{ Cell = "D11"; Value = "0.07734" }
    @{ Cell = "E11"; Value = "  +1.24%  " }
    @{ Cell = "B12"; Value = "Polkadot" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" }
    @{ Cell = "D12"; Value = "4.989" }
    @{ Cell = "E12"; Value = "  -1.04%  " }
    @{ Cell = "B13"; Value = "WrappedEther" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D13"; Value = "1.777.70" }
    @{ Cell = "E13"; Value = "  -3.30%  " }
    @{ Cell = "D14"; Value = "0.6694" }
    @{ Cell = "E14"; Value = "  -1.62%  " }
    @{ Cell = "D15"; Value = "83.03" }
    @{ Cell = "E15"; Value = "  -5.47%  " }
    @{ Cell = "D16"; Value = "6.100" }
    @{ Cell = "E16"; Value = "  -0.29%  " }
    @{ Cell = "D17"; Value = "0.000008395" }
    @{ Cell = "E17"; Value = "  +2.27%  " }
    @{ Cell = "D18"; Value = "29.110.85" }
    @{ Cell = "E18"; Value = "  -0.13%  " }
    @{ Cell = "D19"; Value = "227.29" }
    @{ Cell = "E19"; Value = "  -1.58%  " }
    @{ Cell = "D20"; Value = "12.48" }
    @{ Cell = "E20"; Value = "  -0.29%  " }
    @{ Cell = "E22"; Value = "  -2.14%  " }
    @{ Cell = "E23"; Value = "  -0.01%  " }
    @{ Cell = "D24"; Value = "159.73" }
    @{ Cell = "E24"; Value = "  -0.53%  " }
    @{ Cell = "D25"; Value = "0.1408" }
    @{ Cell = "E25"; Value = "  -1.97%  " }
    @{ Cell = "D26"; Value = "8.630" }
    @{ Cell = "E26"; Value = "  -0.78%  " }
    @{ Cell = "E27"; Value = "  -0.69%  " }
    @{ Cell = "D28"; Value = "1.511" }
    @{ Cell = "E28"; Value = "  +0.58%  " }
    @{ Cell = "E29"; Value = "  -3.54%  " }
    @{ Cell = "D30"; Value = "4.048" }
    @{ Cell = "E30"; Value = "  -2.26%  " }
    @{ Cell = "E31"; Value = "  +0.18%  " }
    @{ Cell = "D32"; Value = "0.05329" }
    @{ Cell = "E32"; Value = "  -0.62%  " }
    @{ Cell = "D33"; Value = "1.876" }
    @{ Cell = "E33"; Value = "  +1.06%  " }
    @{ Cell = "D34"; Value = "0.7573" }
    @{ Cell = "E34"; Value = "  +0.27%  " }
    @{ Cell = "E35"; Value = "  +0.33%  " }
    @{ Cell = "E36"; Value = "  -0.57%  " }
    @{ Cell = "D37"; Value = "1.272.21" }
    @{ Cell = "E38"; Value = "  -1.70%  " }
    @{ Cell = "D39"; Value = "2.734" }
    @{ Cell = "E39"; Value = "  +0.31%  " }
    @{ Cell = "D40"; Value = "0.9282" }
    @{ Cell = "E40"; Value = "  -1.77%  " }
    @{ Cell = "D41"; Value = "0.08935" }
    @{ Cell = "E41"; Value = "  +16.49%  " }
    @{ Cell = "D42"; Value = "5.975" }
    @{ Cell = "E42"; Value = "  -0.52%  " }
    @{ Cell = "D43"; Value = "1.002" }
    @{ Cell = "E43"; Value = "  +0.30%  " }
    @{ Cell = "D44"; Value = "102.79" }
    @{ Cell = "E44"; Value = "  -1.79%  " }
    @{ Cell = "D45"; Value = "1.965.42" }
    @{ Cell = "E45"; Value = "  -1.03%  " }
    @{ Cell = "D46"; Value = "0.5162" }
    @{ Cell = "E46"; Value = "  -0.39%  " }
    @{ Cell = "E47"; Value = "  +0.06%  " }
    @{ Cell = "E48"; Value = "  -1.32%  " }
    @{ Cell = "D49"; Value = "63.32" }
    @{ Cell = "E49"; Value = "  -1.47%  " }
    @{ Cell = "D50"; Value = "0.05912" }
    @{ Cell = "E50"; Value = "  -0.54%  " }
    @{ Cell = "B51"; Value = "EnergySwap" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D51"; Value = "8.835" }
    @{ Cell = "E51"; Value = "  -6.75%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $text = $u.Value
    $looksNumeric = $text -match "^[+-]?\d+\.\d+$"
    if ($looksNumeric) {
        $origStyle = $range.Style
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = $origStyle
    } else {
        $range.Value = $text
    }
}
